# Auto-generated edit script applying numeric updates described by the commit diff.
# Each sheet is addressed by name; values are written cell-by-cell via Range().Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 688.4231
$ws.Range("I18").Value = 722.6087
$ws.Range("J18").Value = 426.33334
$ws.Range("K18").Value = 722.6087
$ws.Range("L18").Value = 426.33334
$ws.Range("M18").Value = -438.6087
$ws.Range("N18").Value = -994.33334
$ws.Range("H69").Value = 50003950
$ws.Range("I69").Value = 2200
$ws.Range("J69").Value = 66671200
$ws.Range("K69").Value = 6600
$ws.Range("L69").Value = 200013600
$ws.Range("M69").Value = -5726
$ws.Range("N69").Value = -200015348
$ws.Range("H72").Value = 50003950
$ws.Range("I72").Value = 2200
$ws.Range("J72").Value = 66671200
$ws.Range("K72").Value = 19800
$ws.Range("L72").Value = 600040800
$ws.Range("M72").Value = -15432
$ws.Range("N72").Value = -600049536
$ws.Range("H115").Value = 7715.6
$ws.Range("I115").Value = 676
$ws.Range("J115").Value = 11506.154
$ws.Range("K115").Value = 2028
$ws.Range("L115").Value = 34518.462
$ws.Range("M115").Value = -461
$ws.Range("N115").Value = -37652.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 21588.75
$ws.Range("J24").Value = 21588.75
$ws.Range("L24").Value = 21588.75
$ws.Range("N24").Value = -22336.75
$ws.Range("H28").Value = 8999.75
$ws.Range("I28").Value = 8999.75
$ws.Range("K28").Value = 8999.75
$ws.Range("M28").Value = -8807.75
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 29993.334
$ws.Range("J94").Value = 29993.334
$ws.Range("L94").Value = 29993.334
$ws.Range("N94").Value = -31795.334
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 14171.5
$ws.Range("J96").Value = 14171.5
$ws.Range("L96").Value = 14171.5
$ws.Range("N96").Value = -19663.5
$ws.Range("H97").Value = 512.9
$ws.Range("I97").Value = 541.125
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 541.125
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = -45.125
$ws.Range("N97").Value = -1392
$ws.Range("H98").Value = 19118.334
$ws.Range("J98").Value = 19118.334
$ws.Range("L98").Value = 19118.334
$ws.Range("N98").Value = -25108.334
$ws.Range("H99").Value = 8999.75
$ws.Range("I99").Value = 8999.75
$ws.Range("K99").Value = 8999.75
$ws.Range("M99").Value = -6004.75
$ws.Range("H100").Value = 21588.75
$ws.Range("J100").Value = 21588.75
$ws.Range("L100").Value = 21588.75
$ws.Range("N100").Value = -23752.75
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H102").Value = 2145.6667
$ws.Range("I102").Value = 1966.6666
$ws.Range("J102").Value = 2503.6667
$ws.Range("K102").Value = 1966.6666
$ws.Range("L102").Value = 2503.6667
$ws.Range("M102").Value = -344.6666
$ws.Range("N102").Value = -5747.6667
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1959.03
$ws.Range("I86").Value = 1966.7291
$ws.Range("J86").Value = 1774.25
$ws.Range("K86").Value = 1966.7291
$ws.Range("L86").Value = 1774.25
$ws.Range("M86").Value = -843.7291
$ws.Range("N86").Value = -4020.25
$ws.Range("H89").Value = 1959.03
$ws.Range("I89").Value = 1966.7291
$ws.Range("J89").Value = 1774.25
$ws.Range("K89").Value = 9833.645500000001
$ws.Range("L89").Value = 8871.25
$ws.Range("M89").Value = -4217.645500000001
$ws.Range("N89").Value = -20103.25
$ws.Range("H94").Value = 2057.375
$ws.Range("I94").Value = 1539.8462
$ws.Range("J94").Value = 4300
$ws.Range("K94").Value = 1539.8462
$ws.Range("L94").Value = 4300
$ws.Range("M94").Value = -1088.8462
$ws.Range("N94").Value = -5202
$ws.Range("H99").Value = 1743.6471
$ws.Range("I99").Value = 1515.8889
$ws.Range("J99").Value = 1999.875
$ws.Range("K99").Value = 1515.8889
$ws.Range("L99").Value = 1999.875
$ws.Range("M99").Value = -17.88889999999992
$ws.Range("N99").Value = -4995.875
$ws.Range("H105").Value = 1734.7646
$ws.Range("I105").Value = 1544.4445
$ws.Range("J105").Value = 1948.875
$ws.Range("K105").Value = 1544.4445
$ws.Range("L105").Value = 1948.875
$ws.Range("M105").Value = 202.5554999999999
$ws.Range("N105").Value = -5442.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 16432.334
$ws.Range("J131").Value = 16432.334
$ws.Range("L131").Value = 16432.334
$ws.Range("N131").Value = -26512.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 6424.0625
$ws.Range("I23").Value = 33.75
$ws.Range("J23").Value = 8554.166999999999
$ws.Range("K23").Value = 101.25
$ws.Range("L23").Value = 25662.501
$ws.Range("M23").Value = 133.75
$ws.Range("N23").Value = -26132.501
$ws.Range("H131").Value = 12195839
$ws.Range("I131").Value = 35714624
$ws.Range("J131").Value = 913.1852
$ws.Range("K131").Value = 107143872
$ws.Range("L131").Value = 2739.5556
$ws.Range("M131").Value = -107138832
$ws.Range("N131").Value = -12819.5556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 19232756
$ws.Range("I97").Value = 2260.7144
$ws.Range("J97").Value = 41668336
$ws.Range("K97").Value = 2260.7144
$ws.Range("L97").Value = 41668336
$ws.Range("M97").Value = -1764.7144
$ws.Range("N97").Value = -41669328

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9830.933999999999
$ws.Range("I93").Value = 2380.7778
$ws.Range("J93").Value = 21006.166
$ws.Range("K93").Value = 2380.7778
$ws.Range("L93").Value = 21006.166
$ws.Range("M93").Value = -1132.7778
$ws.Range("N93").Value = -23502.166
$ws.Range("H100").Value = 19893.967
$ws.Range("I100").Value = 29377.777
$ws.Range("J100").Value = 6762.5386
$ws.Range("K100").Value = 29377.777
$ws.Range("L100").Value = 6762.5386
$ws.Range("M100").Value = -28836.777
$ws.Range("N100").Value = -7844.5386
$ws.Range("H132").Value = 3325775.8
$ws.Range("I132").Value = 3666304
$ws.Range("J132").Value = 5624.5
$ws.Range("K132").Value = 10998912
$ws.Range("L132").Value = 16873.5
$ws.Range("M132").Value = -10996382
$ws.Range("N132").Value = -21933.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10555.182
$ws.Range("I81").Value = 974.5
$ws.Range("J81").Value = 11513.25
$ws.Range("K81").Value = 1949
$ws.Range("L81").Value = 23026.5
$ws.Range("M81").Value = -888
$ws.Range("N81").Value = -25148.5
$ws.Range("H84").Value = 10555.182
$ws.Range("I84").Value = 974.5
$ws.Range("J84").Value = 11513.25
$ws.Range("K84").Value = 9745
$ws.Range("L84").Value = 115132.5
$ws.Range("M84").Value = -4441
$ws.Range("N84").Value = -125740.5

